$d = $word.ActiveDocument

# wdReplaceAll = 2
$wdReplaceAll = 2

function Replace-Text($findText, $replaceText) {
    $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replaceText, $wdReplaceAll)
}

# 1. Title: "PROD (Nov/2025) × STG (Out/2025)" -> "PROD (Out/2025) × PROD (Nov/2025)"
Replace-Text "PROD (Nov/2025) × STG (Out/2025)" "PROD (Out/2025) × PROD (Nov/2025)"

# 2. Visão Geral paragraph: "...ambiente STG em 10/..." -> "...ambiente PROD em 10/..."
Replace-Text "o teste equivalente realizado no ambiente STG em" "o teste equivalente realizado no ambiente PROD em"

# 3. Table header cell: "STG (Out/25)" -> "PROD (Out/25)"
Replace-Text "STG (Out/25)" "PROD (Out/25)"

# 4. "-36 ms (melhor no PROD)" -> "-36 ms (melhor no 11/25)"
Replace-Text "-36 ms (melhor no PROD)" "-36 ms (melhor no 11/25)"

# 5. "-53 ms (melhor no PROD)" -> "-53 ms (melhor no 11/25)"
Replace-Text "-53 ms (melhor no PROD)" "-53 ms (melhor no 11/25)"

# 6. "-920 ms (muito melhor no PROD)" -> "-920 ms (muito melhor no 11/25)"
Replace-Text "-920 ms (muito melhor no PROD)" "-920 ms (muito melhor no 11/25)"

# 7. "picos muito menores no PROD" -> "picos muito menores no mês 10/25"
Replace-Text "picos muito menores no PROD" "picos muito menores no mês 10/25"

# 8. "+0.078 (melhoria significativa no PROD)" -> "+0.078 (melhoria significativa)"
Replace-Text "+0.078 (melhoria significativa no PROD)" "+0.078 (melhoria significativa)"

# 9. Analysis intro paragraph
Replace-Text "A análise comparativa demonstra uma melhoria significativa no ambiente PROD em relação ao teste de estresse realizado no mês anterior no ambiente STG." "A análise comparativa demonstra uma melhoria significativa no ambiente PROD 11/25 em relação ao teste de estresse realizado no mês anterior 10/25."

# 10. Bullet: ambiente PROD consistency
Replace-Text "• O ambiente PROD apresentou tempos muito mais consistentes, especialmente nos " "• O ambiente PROD 11/12 apresentou tempos muito mais consistentes, especialmente nos "

# 11. Bullet: STG picos -> PROD 10/2025
Replace-Text "• Enquanto o STG atingiu picos de até 1,4 segundos no 99th percentile, o PROD manteve todas as respostas abaixo de 0,5 segundos, demonstrando excelente capacidade de processamento." "• Enquanto o PROD 10/2025 atingiu picos de até 1,4 segundos no 99th percentile, o PROD manteve todas as respostas abaixo de 0,5 segundos, demonstrando excelente capacidade de processamento."

# 12. Merge split "throughput" runs (no text change, just a run consolidation)
Replace-Text "• O throughput manteve-se praticamente igual, indicando que a infraestrutura sustenta volume de tráfego estável." "• O throughput manteve-se praticamente igual, indicando que a infraestrutura sustenta volume de tráfego estável."

# 13. Conclusão paragraph
Replace-Text "O Stratega PROD demonstrou desempenho superior ao ambiente STG em todos os indicadores de latência e qualidade percebida pelo usuário." "O Stratega PROD 11/25 demonstrou desempenho superior ao ambiente PROD 10/25 em todos os indicadores de latência e qualidade percebida pelo usuário."

# 14. Checkmark bullet
Replace-Text "✔ O ambiente PROD está mais otimizado e demonstra maior capacidade de escalar." "✔ O ambiente PROD 11/25 está mais otimizado e demonstra maior capacidade de escalar."

# 15. Resultado final
Replace-Text "Resultado final: **PROD apresenta evolução de performance e estabilidade em relação ao mês anterior.**" "Resultado final: **PROD 11/2025 apresenta evolução de performance e estabilidade em relação ao mês anterior.**"
